$d = $word.ActiveDocument

# --- Paragraph with the "m:if" field: turn the field into plain text runs ---
$f1 = $d.Fields.Item(1)
[void]$f1.Delete()
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End
$p2Range = $d.Range($p2Start, $p2End - 1)
$ifXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">{m:if </w:t></w:r><w:r><w:t>self.name}</w:t></w:r><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>The predicate never evaluates to a boolean type ([EClassifier=EString]).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$p2Range.InsertXML($ifXml)

# --- Paragraph with the "m:endif" field: turn the field into plain text ---
$f2 = $d.Fields.Item(1)
[void]$f2.Delete()
$p4 = $d.Paragraphs.Item(4)
$p4Start = $p4.Range.Start
$p4End = $p4.Range.End
$p4Range = $d.Range($p4Start, $p4End - 1)
$endifXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">{m:endif}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$p4Range.InsertXML($endifXml)
